$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two columns that hold R_orig (O) and R_star (R).
# Deleting R first then O avoids re-indexing issues since O is left of R.
$ws.Range("R1").EntireColumn.Delete()
$ws.Range("O1").EntireColumn.Delete()

# Restore the view state seen in the edited workbook (best effort; the
# headless host does not persist all window-scroll metadata).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 11
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("O10").Select()
